# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (currently blank) column
# between the existing "In Advance" (M) and "Late" (N) columns, pushing
# "Late"/"Waived"/"Outstanding" one column to the right. The sheet also
# becomes the active/selected sheet in the workbook (it previously was
# "NewLoanInput").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column immediately to the left (M, "In
# Advance") so the freshly inserted column inherits it, matching what
# Excel does when a new column is inserted next to existing data.
$leftColWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new blank column at N; everything from N onward (Late, Waived,
# Outstanding) shifts one column to the right.
$ws.Columns("N:N").Insert()

# Give the newly inserted column the same width as its left neighbor.
$ws.Columns("N:N").ColumnWidth = $leftColWidth

# Make "Repayment schedule" the active sheet/tab, and move the selection
# to what is now the bottom-right data cell (old selection G9 on the
# "Outstanding" column, now shifted to S9... i.e. the cell the user ended
# up on after inserting the column).
$ws.Activate()
$ws.Range("S9").Select()
